# Regenerate the "K" (strikeouts) column (column G) values for the
# chirinos_yonny save_data sheet, replacing the old Strike#-derived values
# with the freshly-calculated K values (and other regen'd stats).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (header "K"), rows 2-26
$newK = @{
    2  = 5
    3  = 5
    4  = 5
    5  = 5
    6  = 3
    7  = 1
    8  = 2
    9  = 1
    10 = 4
    11 = 3
    12 = 3
    13 = 2
    14 = 4
    15 = 2
    16 = 1
    17 = 0
    18 = 3
    19 = 2
    20 = 1
    21 = 2
    22 = 1
    23 = 4
    24 = 1
    25 = 2
    26 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
